# RE-L00-Organization.pptx — "squash! squash! Add New Lecture Plan"
#
# The lecture-plan table on the course-schedule slide lists the exercise
# sheets handed out each week. The row for the Elicitation exercises
# mistakenly names both the first AND the second sheet "E02"; this fixes
# the second one to "E03" (Elicitation I stays E02, Elicitation II becomes
# E03), matching the numbering used for every other exercise in the table.

$p = $ppt.ActivePresentation

# Find the slide that contains the lecture-plan table (the one with the
# "Publication Date" / "Submission Deadline" / "Exercise" header row) so the
# script does not depend on a hard-coded slide index.
$targetCell = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cell = $tbl.Cell($r, $c)
                    $txt = $cell.Shape.TextFrame.TextRange.Text
                    if ($txt -like "*E02*Elicitation I*E02*Elicitation II*") {
                        $targetCell = $cell
                    }
                }
            }
        }
    }
}

if ($targetCell -ne $null) {
    $tr = $targetCell.Shape.TextFrame.TextRange
    $tr.Text = "E02 – Elicitation I, E03 – Elicitation II"
}
